$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Columns("D:D").Insert()
$ws.Range("D1").Value = "Thời vụ"
$ws.Range("D1").Select()
